$d = $word.ActiveDocument

# --- Change 1: expand the commentary on the excluded road_density variable ---
$find = $d.Content.Find
$find.Execute(
    "n’apparaît pas dans les modèles testés sans justification. Cette exclusion pourrait biaiser les résultats.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "n’apparaît pas dans les modèles testés sans justification. Cette exclusion est préoccupante car la densité routière est un indicateur clé du stress urbain pour les oiseaux, affectant la pollution sonore, la qualité de l’air, et le risque de collision (Kight & Swaddle, 2011). Son omission pourrait biaiser les estimations des autres coefficients.",
    2
) | Out-Null

# --- Change 2: add the Kight & Swaddle (2011) reference to the bibliography ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Hartig, F.*DHARMa*") {
        $p.Range.InsertParagraphAfter() | Out-Null
        $newPara = $p.Next()
        $npr = $newPara.Range

        $xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1007"/>
    </w:numPr>
    <w:pStyle w:val="Compact"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Kight, C. R., &amp; Swaddle, J. P. (2011). How and why environmental noise impacts animals: an integrative, mechanistic review.</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:iCs/>
      <w:i/>
    </w:rPr>
    <w:t xml:space="preserve">Ecology Letters</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">, 14(10), 1052-1061.</w:t>
  </w:r>
</w:p>
"@

        $npr.InsertXML($xml) | Out-Null
        break
    }
}
